# Add data for 2022-04-16
# Rename the "Through" sheet/title/header to reflect the new as-of date
# and bump the counts for the neighborhoods that had a carjacking on 2022-04-16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab and its header label (column B, row 1) plus the
# shared "April 2022 (through April 07)" string everywhere it appears.
$ws.Name = "Through 2022-04-08"
$ws.Range("B1").Value = "April 2022 (through April 08)"

function Add-Count {
    param(
        [string]$CellRef,
        [int]$Delta
    )
    $cell = $ws.Range($CellRef)
    $current = $cell.Value
    if ($null -eq $current) {
        $current = 0
    }
    $cell.Value = $current + $Delta
}

# Column B = "April 2022 (through April 08)" -- new incident's current month
Add-Count "B6" 1   # Humboldt Park
Add-Count "B9" 1   # Loop
Add-Count "B23" 1  # Auburn Gresham

# Historical "April" columns for the same neighborhoods, updated for
# year-over-year comparison alongside the new data point.
Add-Count "N2" 1   # Austin - April 2019
Add-Count "R2" 1   # Austin - April 2018
Add-Count "AD3" 1  # Englewood - April 2015
Add-Count "V4" 1   # North Lawndale - April 2017
Add-Count "V5" 1   # Garfield Park - April 2017
Add-Count "N8" 1   # Chicago Lawn - April 2019
Add-Count "J11" 1  # Chatham - April 2020
Add-Count "R17" 1  # Belmont Cragin - April 2018
Add-Count "V17" 1  # Belmont Cragin - April 2017
Add-Count "F20" 1  # Near South Side - April 2021
Add-Count "F34" 1  # South Deering - April 2021
Add-Count "J86" 1  # South Chicago - April 2020
Add-Count "Z89" 1  # United Center - April 2016
